# Commit: update monthly and yearly outcome_ath dictionaries
# - Fix "identifer" -> "identifier" typo in two variable descriptions
#   on the "Variables" sheet.
# - Update the sheet's saved selection to active cell D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("D2").Value = "Unique identifier for the row in Opal"
$ws.Range("D3").Value = "Unique identifier for the child"

$ws.Activate()
$ws.Range("D4").Select() | Out-Null
